# Update the "fixed" date placeholder text (the a:fld/datetimeFigureOut
# field's cached display text) from 1/24/21 to 6/18/25 everywhere it
# appears: the slide master, every slide layout, the handout master and
# the notes master.

$p = $ppt.ActivePresentation
$oldDate = "1/24/21"
$newDate = "6/18/25"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster

# Every slide layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j)
}

# Handout master.
if ($p.HasHandoutMaster) {
    Update-DatePlaceholder $p.HandoutMaster
}

# Notes master.
if ($p.HasNotesMaster) {
    Update-DatePlaceholder $p.NotesMaster
}
